$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: bit-size headers changed from 16비트/16비트/32비트 to 8비트/8비트/48비트
$ws.Range("A1").Value = "8비트"
$ws.Range("B1").Value = "8비트"
$ws.Range("C1").Value = "48비트"

# Row 2: field names stay the same (WorldID, ServerID, value)
$ws.Range("A2").Value = "WorldID"
$ws.Range("B2").Value = "ServerID"
$ws.Range("C2").Value = "value"

# Row 4: description text updated with new ID capacity
$ws.Range("A4").Value = "각 서버별 281474976710655 까지 ID 발급 가능 (약 280조)"

# Update the active selection to H7 as in the target workbook
$ws.Range("H7").Select()

$wb.Save()
